$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe forces Excel to store the value as literal text
# (prevents numeric-looking strings like '578.48' or '13.10' from being
# auto-converted to numbers and losing formatting such as trailing zeros).
$apos = [string][char]0x0027
$u2083 = [string][char]0x2083

$ws.Range("D2").Value = [string]$apos + '60.522.15'
$ws.Range("E2").Value = [string]'  -5.47%  '
$ws.Range("D3").Value = [string]$apos + '3.007.62'
$ws.Range("E3").Value = [string]'  -6.84%  '
$ws.Range("E4").Value = [string]'  -0.01%  '
$ws.Range("D5").Value = [string]$apos + '578.48'
$ws.Range("E5").Value = [string]'  -2.87%  '
$ws.Range("D6").Value = [string]$apos + '126.66'
$ws.Range("E6").Value = [string]'  -8.17%  '
$ws.Range("E7").Value = [string]'  +0.14%  '
$ws.Range("D8").Value = [string]$apos + '3.000.15'
$ws.Range("E8").Value = [string]'  -6.93%  '
$ws.Range("E9").Value = [string]'  -3.34%  '
$ws.Range("E10").Value = [string]'  -7.96%  '
$ws.Range("D11").Value = [string]$apos + '5.15'
$ws.Range("E11").Value = [string]'  -3.49%  '
$ws.Range("E12").Value = [string]'  -3.86%  '
$ws.Range("E13").Value = [string]'  -7.70%  '
$ws.Range("D14").Value = [string]$apos + '32.65'
$ws.Range("E14").Value = [string]'  -7.24%  '
$ws.Range("E15").Value = [string]'  +0.19%  '
$ws.Range("D16").Value = [string]$apos + '3.501.75'
$ws.Range("E16").Value = [string]'  -6.90%  '
$ws.Range("D17").Value = [string]$apos + '3.017.78'
$ws.Range("E17").Value = [string]'  -6.52%  '
$ws.Range("D18").Value = [string]$apos + '60.479.86'
$ws.Range("E18").Value = [string]'  -5.57%  '
$ws.Range("E19").Value = [string]'  -4.07%  '
$ws.Range("D20").Value = [string]$apos + '431.98'
$ws.Range("E20").Value = [string]'  -7.73%  '
$ws.Range("D21").Value = [string]$apos + '13.10'
$ws.Range("E21").Value = [string]'  -7.30%  '
$ws.Range("D22").Value = [string]$apos + '0.664'
$ws.Range("E22").Value = [string]'  -6.31%  '
$ws.Range("D23").Value = [string]$apos + '7.03'
$ws.Range("E23").Value = [string]'  -9.34%  '
$ws.Range("D24").Value = [string]$apos + '12.81'
$ws.Range("E24").Value = [string]'  -5.31%  '
$ws.Range("D25").Value = [string]$apos + '79.52'
$ws.Range("E25").Value = [string]'  -5.01%  '
$ws.Range("E26").Value = [string]'  -0.03%  '
$ws.Range("E27").Value = [string]'  -0.18%  '
$ws.Range("D28").Value = [string]$apos + '2.57'
$ws.Range("E28").Value = [string]'  -4.98%  '
$ws.Range("D29").Value = [string]$apos + '7.34'
$ws.Range("E29").Value = [string]'  -7.07%  '
$ws.Range("E30").Value = [string]'  -8.56%  '
$ws.Range("D31").Value = [string]$apos + '6.16'
$ws.Range("E31").Value = [string]'  -10.50%  '
$ws.Range("D32").Value = [string]$apos + '25.32'
$ws.Range("E32").Value = [string]'  -8.69%  '
$ws.Range("D33").Value = [string]$apos + '0.0934'
$ws.Range("E33").Value = [string]'  -9.84%  '
$ws.Range("D34").Value = [string]$apos + '2.15'
$ws.Range("E34").Value = [string]'  -12.36%  '
$ws.Range("E35").Value = [string]'  -8.42%  '
$ws.Range("E36").Value = [string]'  -5.77%  '
$ws.Range("D37").Value = [string]$apos + '50.01'
$ws.Range("E37").Value = [string]'  -3.37%  '
$ws.Range("D38").Value = [string]$apos + '0.0' + $u2083 + '0664'
$ws.Range("E38").Value = [string]'  -9.79%  '
$ws.Range("D39").Value = [string]$apos + '8.35'
$ws.Range("E39").Value = [string]'  +2.21%  '
$ws.Range("D40").Value = [string]$apos + '0.0360'
$ws.Range("E40").Value = [string]'  -8.93%  '
$ws.Range("B41").Value = [string]'Kaspa'
$ws.Range("C41").Value = [string]'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D41").Value = [string]$apos + '0.110'
$ws.Range("E41").Value = [string]'  -3.03%  '
$ws.Range("B42").Value = [string]'Bittensor'
$ws.Range("C42").Value = [string]'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D42").Value = [string]$apos + '386.96'
$ws.Range("E42").Value = [string]'  -5.09%  '
$ws.Range("D43").Value = [string]$apos + '2.49'
$ws.Range("E43").Value = [string]'  -10.68%  '
$ws.Range("D44").Value = [string]$apos + '2.658.20'
$ws.Range("E44").Value = [string]'  -6.96%  '
$ws.Range("E46").Value = [string]'  -8.24%  '
$ws.Range("D47").Value = [string]$apos + '2.01'
$ws.Range("E47").Value = [string]'  -7.45%  '
$ws.Range("D48").Value = [string]$apos + '118.02'
$ws.Range("E48").Value = [string]'  -7.55%  '
$ws.Range("E49").Value = [string]'  -4.64%  '
$ws.Range("D50").Value = [string]$apos + '23.75'
$ws.Range("E50").Value = [string]'  -8.30%  '
$ws.Range("D51").Value = [string]$apos + '0.134'
$ws.Range("E51").Value = [string]'  +2.05%  '
